$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on Price/Volume columns so numeric-looking strings
# (e.g. "1.01") are not auto-coerced into actual numbers by Excel, matching
# the workbook's original inline-string cell typing.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.399.29"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "2.223.07"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "111.23"
$ws.Range("E5").Value = "  -6.68%  "
$ws.Range("D6").Value = "290.00"
$ws.Range("E6").Value = "  +8.00%  "
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -2.24%  "
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").Value = "0.599"
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("D10").Value = "43.82"
$ws.Range("E10").Value = "  -7.07%  "
$ws.Range("D11").Value = "0.0912"
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("D12").Value = "54.31"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "8.64"
$ws.Range("E13").Value = "  -7.16%  "
$ws.Range("D14").Value = "1.02"
$ws.Range("E14").Value = "  +12.74%  "
$ws.Range("E15").Value = "  -2.73%  "
$ws.Range("D16").Value = "14.91"
$ws.Range("E16").Value = "  -4.55%  "
$ws.Range("D17").Value = "2.560.83"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("D18").Value = "2.262.54"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "42.452.58"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("E20").Value = "  -3.38%  "
$ws.Range("D21").Value = "7.14"
$ws.Range("E21").Value = "  +4.67%  "
$ws.Range("D22").Value = "73.02"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("D23").Value = "3.33"
$ws.Range("E23").Value = "  +12.85%  "
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").Value = "235.13"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "9.00"
$ws.Range("E26").Value = "  -7.24%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").Value = "11.42"
$ws.Range("E28").Value = "  -7.29%  "
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("D30").Value = "37.62"
$ws.Range("E30").Value = "  -10.25%  "
$ws.Range("D31").Value = "173.27"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").Value = "3.10"
$ws.Range("E32").Value = "  -7.21%  "
$ws.Range("D33").Value = "21.20"
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("D34").Value = "0.0879"
$ws.Range("E34").Value = "  -3.44%  "
$ws.Range("D35").Value = "5.62"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").Value = "4.99"
$ws.Range("E36").Value = "  +6.51%  "
$ws.Range("E37").Value = "  -5.46%  "
$ws.Range("E38").Value = "  -3.10%  "
$ws.Range("D39").Value = "0.0378"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").Value = "0.105"
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("D41").Value = "2.39"
$ws.Range("E41").Value = "  -5.84%  "
$ws.Range("D42").Value = "71.82"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "0.231"
$ws.Range("E43").Value = "  -4.14%  "
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "12.38"
$ws.Range("E45").Value = "  -9.73%  "
$ws.Range("E46").Value = "  -3.49%  "
$ws.Range("D47").Value = "5.35"
$ws.Range("E47").Value = "  -6.15%  "
$ws.Range("D48").Value = "1.27"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "1.67"
$ws.Range("E49").Value = "  +3.03%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "8.40"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("D51").Value = "100.64"
$ws.Range("E51").Value = "  -1.26%  "

# Restore the default cell style (removes the temporary text-format style
# index from the cells so only Value content actually differs).
$ws.Range("D2:E51").Style = "Normal"

